$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column P (16th column); this shifts the
# existing P:Y columns (and their data / widths) one slot to the right,
# becoming Q:Z.
$ws.Columns.Item(16).Insert()

# Give the freshly inserted column P the same width the old column Q had
# (closest value reachable through the ColumnWidth property).
$ws.Columns.Item(16).ColumnWidth = 4.571428571428571

# Populate the new column's header (row 1) and data (row 2).
$ws.Cells.Item(1, 16).Value = "fgfg"
$ws.Cells.Item(2, 16).Value = 0
